$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Variant A" / "Winner (Express)" columns for both data rows now all read
# "With express lanes" (the express-lanes variant won the test).
$ws.Range("I2").Value = "With express lanes"
$ws.Range("J2").Value = "With express lanes"
$ws.Range("K2").Value = "With express lanes"

$ws.Range("I3").Value = "With express lanes"
$ws.Range("J3").Value = "With express lanes"
$ws.Range("K3").Value = "With express lanes"

# Move the active selection to I9, scrolled so column B is leftmost.
$ws.Range("I9").Select()
